$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 14 new rows before the old row 69 (the 5 blank rows + total row get
# pushed down to 82-88; the old total-row formula/style/number format travels
# down automatically with the insert). ---
$ws.Rows("69:82").Insert() | Out-Null

# --- Column A width (user resized the column; also drops the old bestFit flag) ---
$ws.Columns("A").ColumnWidth = 11.6

# --- New expense/work-log rows 69-81 -------------------------------------
# Each row: Date (col A), Hours (col C), Price (col D), Description (col E).
# The shared-string table append order follows the order the text was typed,
# not the row order, so the E-column assignments below are intentionally
# sequenced to match that (69,70,71,72,73,79,80,76,77,75,74,78,81).

function Set-Row($r, $date, $hours, $price) {
    $ws.Cells.Item($r, 1).Value = $date
    $ws.Cells.Item($r, 3).Value = $hours
    $ws.Cells.Item($r, 4).Value = $price
}

Set-Row 69 42443 3 150
Set-Row 70 42443 0.5 25
Set-Row 71 42444 2 50
Set-Row 72 42444 0.5 25
Set-Row 73 42444 1.5 75
Set-Row 74 42445 4 200
Set-Row 75 42446 3 150
Set-Row 76 42447 0.5 25
Set-Row 77 42447 1.5 75
Set-Row 78 42447 0.5 25
Set-Row 79 42448 2 100
Set-Row 80 42448 2 100
Set-Row 81 42450 4 200

$ws.Range("E69").Value = "Реализиран интерфейс за отпуски на служител, и полагаем годишен отпуск."
$ws.Range("E70").Value = "Добавени годишни отпуски за всички служители."
$ws.Range("E71").Value = "Проектиране на визуализациите на видовете графици, номенклатура на видовете графици и генерация на графиците"
$ws.Range("E72").Value = "Реализация на номенклатура за видове графици"
$ws.Range("E73").Value = "Реализация на визуализацията на екипи съчетани с графици, филтър по дата и по вид график"
$ws.Range("E79").Value = "Реализация на нова логика за бързо визуализиране на графиците и екипите в организацията"
$ws.Range("E80").Value = "Реализация на нов потребителски интерфейс за бързо визуализиране на графиците и екипите в организацията."
$ws.Range("E76").Value = "Добавени полетата ТРЗ код и Графици код"
$ws.Range("E77").Value = "Довършено генерирането на графици с правилна формула"
$ws.Range("E75").Value = "Реализация на потребителски интерфейс за генериране на графиците. Тестване на логиката на генерираенто, чистене на проблеми, разписана логика за запаметяване на графици."
$ws.Range("E74").Value = "Реализация на логиката за генериране на графици, вмъкване на отпуските в графиците"
$ws.Range("E78").Value = "Добавено описание на редовете въф формата за работни дни през месеца."
$ws.Range("E81").Value = "Генерация на графици."

# Row heights for the wrapped-text description cells (best-effort match of
# Excel's auto-fit wrap height for the new rows).
$ws.Rows(69).RowHeight = 30
$ws.Rows(71).RowHeight = 45
$ws.Rows(73).RowHeight = 30
$ws.Rows(74).RowHeight = 30
$ws.Rows(75).RowHeight = 60
$ws.Rows(77).RowHeight = 30
$ws.Rows(78).RowHeight = 30
$ws.Rows(79).RowHeight = 45
$ws.Rows(80).RowHeight = 45

# --- Fix the grand-total formula, now living on row 88, to cover the newly
# inserted data rows (D2:D81 instead of D2:D68). ---
$ws.Range("D88").Formula = "=SUM(D2:D81)"

# --- Sheet view: scrolled down, zoomed in, new selection ---
$excel.ActiveWindow.Zoom = 115
$excel.ActiveWindow.ScrollRow = 72
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D82").Select() | Out-Null
